# Fix a typo ("representaividad" -> "representatividad") and add a new
# paragraph right after it in the body placeholder of slide 55.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(55)
$shape = $s.Shapes.Item("Content Placeholder 2")

$cr = [char]13

$tr = $shape.TextFrame.TextRange
$newText = "El valor de R-squared es subjetivo dependiendo del tope o techo o criterio que le de el investigador a la representatividad de x con respecto a y en el modelo de regresión lineal." + $cr + `
    "Significa que tanto explica la variable independiente x a la variable dependiente y" + $cr + `
    "Si el valor fue del 70% entonces el modelo está por debajo de las expectativas." + $cr + `
    "Si el valor fue del 60% entonces el modelo está dentro de las expectativas esperadas."

$tr.Text = $newText
